# Update crypto price/volume/coin data per latest scrape (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.171.37'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.785.00'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.09'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.547'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.13'
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0688'
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.042.16'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.97'
$ws.Range('E13').Value = '  -4.13%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.736.73'
$ws.Range('E14').Value = '  -3.06%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.624'
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.160.01'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.67'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.73'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.12'
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.01'
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.15'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  +1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0520'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('E33').Value = '  +3.26%  '
$ws.Range('E34').Value = '  -1.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.448.50'
$ws.Range('E35').Value = '  +3.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.56'
$ws.Range('E36').Value = '  +10.25%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '83.88'
$ws.Range('E40').Value = '  +5.00%  '
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.65'
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('E45').Value = '  +2.30%  '
$ws.Range('E46').Value = '  +0.47%  '
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.941.61'
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('E49').Value = '  -4.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.95'
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('E51').Value = '  +0.20%  '
